# Update presigned Samsara S3 URLs in column K (video_Interior) for rows 2-6.
# Only the X-Amz-Date, response-expires timestamp, and X-Amz-Signature query
# parameters change (the pre-signed URL was regenerated at a later time);
# everything else in each URL stays identical.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276206556/TJaAJ2AJpb-camera-video-segment-driver-1743276209056.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSLXPPTF5L%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T100006Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJGMEQCIByNgyG0wrR2KKpwVEVxp0qNKwd2eE1LZUo8AILnLcw9AiAojFkOOtGUz4MQ19qPA%2Fiat3117fI%2B7c8lVt25pNXizCrmAwiF%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM0frgMjiSYxgjpr7HKroDa7cdZ4E00FU%2FGVpS4encjtSDCUiIpLlhGUPqALd30WCtTpSOqOTrkSxzugUTh7SBcj1Lk8zBSD6KQ8HLxXrMJhKRDXkjrZCJBbs7fSIjDTNChpSfripixTeZ7t4D2LAN%2F9aAoVsUXx3ttPk4X9qTtahDL8uZlJYAtHn3XjGEaxUvV5liPV76Zt5g5cavEUkjHpF4ANgWEbmBAR8J45WxZg6cEYjJiiUxFATYGKXR%2FqP9RFNxSOOj%2FWvwSSyg6rjSw1bJGTjr%2FwkDjNTaRthUfQpRrGUm4a%2B6ezscUhjuTXTEaGD7OPR2j2bpFzznQP4w9NRIMHcFfDnWjro2SFN7FYWIi8VkaPc02wdf8Omv%2FsrZKVnNGLMbHSVKJDop6ERLBURHCvyeCNEJ7Kh0kA6dROAABJSVGILYTtZD74nlnDDNsJu3nx2iscpEm%2FfKwgwtP8e6M0YSFofa%2B201q%2BKf4D9MPiJ2PT0IWnRdFNyMtHYFufI4TIE3YOmo6ALsk9CtdgiUvWCtuiZ82hZmtudLYc6uv9gkYFi84D6A8DRXKU1hY9zzIKuJREYR5yFz7isct0%2FAg6a3RYMAOzC%2FiaO%2FBjqmATOYI%2FYR%2B9o46ZcKxgoMQPdKP8B9xS61DGgRun49TZ3lihG%2BKl8pyggLbpYcYW2gvGPDIVnCsbrtG%2FgDBQB7b60ItGso5LeAzOiTvvRc8y%2Fl667zX3oCmBNpLYhdQerjD4WJXwCRu%2BR6009VTYMEYQ9GYTHUePdfzsmvrZ8d%2B%2BUWS7nDZl%2F40%2Bm6ufWvB8l2yGAaSOZhB%2FYI3VZ6Z4sOrO1i5%2BEnpFI%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2018%3A00%3A06%20GMT&X-Amz-Signature=1de742aa0b7cbea2f3bcf8180c6a31107832ad7cf604fb1fc63dd9255bc7e3c2"

$ws.Range("K3").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276164863/OF6jgx1Lcp-camera-video-segment-driver-1743276167363.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSLXPPTF5L%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T100006Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJGMEQCIByNgyG0wrR2KKpwVEVxp0qNKwd2eE1LZUo8AILnLcw9AiAojFkOOtGUz4MQ19qPA%2Fiat3117fI%2B7c8lVt25pNXizCrmAwiF%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM0frgMjiSYxgjpr7HKroDa7cdZ4E00FU%2FGVpS4encjtSDCUiIpLlhGUPqALd30WCtTpSOqOTrkSxzugUTh7SBcj1Lk8zBSD6KQ8HLxXrMJhKRDXkjrZCJBbs7fSIjDTNChpSfripixTeZ7t4D2LAN%2F9aAoVsUXx3ttPk4X9qTtahDL8uZlJYAtHn3XjGEaxUvV5liPV76Zt5g5cavEUkjHpF4ANgWEbmBAR8J45WxZg6cEYjJiiUxFATYGKXR%2FqP9RFNxSOOj%2FWvwSSyg6rjSw1bJGTjr%2FwkDjNTaRthUfQpRrGUm4a%2B6ezscUhjuTXTEaGD7OPR2j2bpFzznQP4w9NRIMHcFfDnWjro2SFN7FYWIi8VkaPc02wdf8Omv%2FsrZKVnNGLMbHSVKJDop6ERLBURHCvyeCNEJ7Kh0kA6dROAABJSVGILYTtZD74nlnDDNsJu3nx2iscpEm%2FfKwgwtP8e6M0YSFofa%2B201q%2BKf4D9MPiJ2PT0IWnRdFNyMtHYFufI4TIE3YOmo6ALsk9CtdgiUvWCtuiZ82hZmtudLYc6uv9gkYFi84D6A8DRXKU1hY9zzIKuJREYR5yFz7isct0%2FAg6a3RYMAOzC%2FiaO%2FBjqmATOYI%2FYR%2B9o46ZcKxgoMQPdKP8B9xS61DGgRun49TZ3lihG%2BKl8pyggLbpYcYW2gvGPDIVnCsbrtG%2FgDBQB7b60ItGso5LeAzOiTvvRc8y%2Fl667zX3oCmBNpLYhdQerjD4WJXwCRu%2BR6009VTYMEYQ9GYTHUePdfzsmvrZ8d%2B%2BUWS7nDZl%2F40%2Bm6ufWvB8l2yGAaSOZhB%2FYI3VZ6Z4sOrO1i5%2BEnpFI%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2018%3A00%3A06%20GMT&X-Amz-Signature=c20edbd04d4da79ff8580f5826e219df7e7bd8b486b77eca77a1c09970b2b11a"

$ws.Range("K4").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276108738/gAxSKLP636-camera-video-segment-driver-1743276111238.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSLXPPTF5L%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T100006Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJGMEQCIByNgyG0wrR2KKpwVEVxp0qNKwd2eE1LZUo8AILnLcw9AiAojFkOOtGUz4MQ19qPA%2Fiat3117fI%2B7c8lVt25pNXizCrmAwiF%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM0frgMjiSYxgjpr7HKroDa7cdZ4E00FU%2FGVpS4encjtSDCUiIpLlhGUPqALd30WCtTpSOqOTrkSxzugUTh7SBcj1Lk8zBSD6KQ8HLxXrMJhKRDXkjrZCJBbs7fSIjDTNChpSfripixTeZ7t4D2LAN%2F9aAoVsUXx3ttPk4X9qTtahDL8uZlJYAtHn3XjGEaxUvV5liPV76Zt5g5cavEUkjHpF4ANgWEbmBAR8J45WxZg6cEYjJiiUxFATYGKXR%2FqP9RFNxSOOj%2FWvwSSyg6rjSw1bJGTjr%2FwkDjNTaRthUfQpRrGUm4a%2B6ezscUhjuTXTEaGD7OPR2j2bpFzznQP4w9NRIMHcFfDnWjro2SFN7FYWIi8VkaPc02wdf8Omv%2FsrZKVnNGLMbHSVKJDop6ERLBURHCvyeCNEJ7Kh0kA6dROAABJSVGILYTtZD74nlnDDNsJu3nx2iscpEm%2FfKwgwtP8e6M0YSFofa%2B201q%2BKf4D9MPiJ2PT0IWnRdFNyMtHYFufI4TIE3YOmo6ALsk9CtdgiUvWCtuiZ82hZmtudLYc6uv9gkYFi84D6A8DRXKU1hY9zzIKuJREYR5yFz7isct0%2FAg6a3RYMAOzC%2FiaO%2FBjqmATOYI%2FYR%2B9o46ZcKxgoMQPdKP8B9xS61DGgRun49TZ3lihG%2BKl8pyggLbpYcYW2gvGPDIVnCsbrtG%2FgDBQB7b60ItGso5LeAzOiTvvRc8y%2Fl667zX3oCmBNpLYhdQerjD4WJXwCRu%2BR6009VTYMEYQ9GYTHUePdfzsmvrZ8d%2B%2BUWS7nDZl%2F40%2Bm6ufWvB8l2yGAaSOZhB%2FYI3VZ6Z4sOrO1i5%2BEnpFI%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2018%3A00%3A06%20GMT&X-Amz-Signature=cf02f4e71397d6aa695129d56ae1ddcd48f71151169d051c17c9ec09f850cc39"

$ws.Range("K5").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205341/1743265652740/hq4516Mzws-camera-video-segment-driver-1743265655240.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSLXPPTF5L%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T100006Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJGMEQCIByNgyG0wrR2KKpwVEVxp0qNKwd2eE1LZUo8AILnLcw9AiAojFkOOtGUz4MQ19qPA%2Fiat3117fI%2B7c8lVt25pNXizCrmAwiF%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM0frgMjiSYxgjpr7HKroDa7cdZ4E00FU%2FGVpS4encjtSDCUiIpLlhGUPqALd30WCtTpSOqOTrkSxzugUTh7SBcj1Lk8zBSD6KQ8HLxXrMJhKRDXkjrZCJBbs7fSIjDTNChpSfripixTeZ7t4D2LAN%2F9aAoVsUXx3ttPk4X9qTtahDL8uZlJYAtHn3XjGEaxUvV5liPV76Zt5g5cavEUkjHpF4ANgWEbmBAR8J45WxZg6cEYjJiiUxFATYGKXR%2FqP9RFNxSOOj%2FWvwSSyg6rjSw1bJGTjr%2FwkDjNTaRthUfQpRrGUm4a%2B6ezscUhjuTXTEaGD7OPR2j2bpFzznQP4w9NRIMHcFfDnWjro2SFN7FYWIi8VkaPc02wdf8Omv%2FsrZKVnNGLMbHSVKJDop6ERLBURHCvyeCNEJ7Kh0kA6dROAABJSVGILYTtZD74nlnDDNsJu3nx2iscpEm%2FfKwgwtP8e6M0YSFofa%2B201q%2BKf4D9MPiJ2PT0IWnRdFNyMtHYFufI4TIE3YOmo6ALsk9CtdgiUvWCtuiZ82hZmtudLYc6uv9gkYFi84D6A8DRXKU1hY9zzIKuJREYR5yFz7isct0%2FAg6a3RYMAOzC%2FiaO%2FBjqmATOYI%2FYR%2B9o46ZcKxgoMQPdKP8B9xS61DGgRun49TZ3lihG%2BKl8pyggLbpYcYW2gvGPDIVnCsbrtG%2FgDBQB7b60ItGso5LeAzOiTvvRc8y%2Fl667zX3oCmBNpLYhdQerjD4WJXwCRu%2BR6009VTYMEYQ9GYTHUePdfzsmvrZ8d%2B%2BUWS7nDZl%2F40%2Bm6ufWvB8l2yGAaSOZhB%2FYI3VZ6Z4sOrO1i5%2BEnpFI%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2018%3A00%3A06%20GMT&X-Amz-Signature=44d9cd5cdd6f1fb184f36ab19a51fe22aa95f4e7a8e5d38f38916df7c5c30856"

$ws.Range("K6").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743255999969/Yjq7ZuVo4p-camera-video-segment-driver-1743256002469.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSLXPPTF5L%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T100006Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJGMEQCIByNgyG0wrR2KKpwVEVxp0qNKwd2eE1LZUo8AILnLcw9AiAojFkOOtGUz4MQ19qPA%2Fiat3117fI%2B7c8lVt25pNXizCrmAwiF%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM0frgMjiSYxgjpr7HKroDa7cdZ4E00FU%2FGVpS4encjtSDCUiIpLlhGUPqALd30WCtTpSOqOTrkSxzugUTh7SBcj1Lk8zBSD6KQ8HLxXrMJhKRDXkjrZCJBbs7fSIjDTNChpSfripixTeZ7t4D2LAN%2F9aAoVsUXx3ttPk4X9qTtahDL8uZlJYAtHn3XjGEaxUvV5liPV76Zt5g5cavEUkjHpF4ANgWEbmBAR8J45WxZg6cEYjJiiUxFATYGKXR%2FqP9RFNxSOOj%2FWvwSSyg6rjSw1bJGTjr%2FwkDjNTaRthUfQpRrGUm4a%2B6ezscUhjuTXTEaGD7OPR2j2bpFzznQP4w9NRIMHcFfDnWjro2SFN7FYWIi8VkaPc02wdf8Omv%2FsrZKVnNGLMbHSVKJDop6ERLBURHCvyeCNEJ7Kh0kA6dROAABJSVGILYTtZD74nlnDDNsJu3nx2iscpEm%2FfKwgwtP8e6M0YSFofa%2B201q%2BKf4D9MPiJ2PT0IWnRdFNyMtHYFufI4TIE3YOmo6ALsk9CtdgiUvWCtuiZ82hZmtudLYc6uv9gkYFi84D6A8DRXKU1hY9zzIKuJREYR5yFz7isct0%2FAg6a3RYMAOzC%2FiaO%2FBjqmATOYI%2FYR%2B9o46ZcKxgoMQPdKP8B9xS61DGgRun49TZ3lihG%2BKl8pyggLbpYcYW2gvGPDIVnCsbrtG%2FgDBQB7b60ItGso5LeAzOiTvvRc8y%2Fl667zX3oCmBNpLYhdQerjD4WJXwCRu%2BR6009VTYMEYQ9GYTHUePdfzsmvrZ8d%2B%2BUWS7nDZl%2F40%2Bm6ufWvB8l2yGAaSOZhB%2FYI3VZ6Z4sOrO1i5%2BEnpFI%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2018%3A00%3A06%20GMT&X-Amz-Signature=d2516e58bb6f9d41d53ac4facdaeed37e571c3b1b0519217c96ac9ca7fe9b86e"

$wb.Save()
